$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")
Write-Host $ws.Name
$cell = $ws.Cells.Item(18,3)
Write-Host ("Text: " + $cell.Text)
Write-Host ("Value2: " + $cell.Value2)
$cell.Value2 = "HelloTest"
Write-Host ("Value2 after set: " + $cell.Value2)
